$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF ("Date") held the literal text "5-13-2011-12" for every data
# row (2-31). The NBA stats for this date were actually pulled one day off
# because of how the box scores were reported, so the value is corrected to
# the proper date text "2012-05-13".
$dateRange = $ws.Range("BF2:BF31")

# Mark the range as Text first so the engine doesn't reinterpret the
# date-shaped literal as a real date serial number when it's assigned.
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 58).Value = "2012-05-13"
}

# Put the cells' style back to the workbook default so no stray
# number-format/style lingers on these cells (they were unstyled before).
$dateRange.Style = "Normal"
